# Apply the VSIG Trial Balance "second page" update:
#  - Populate the report header block (E2:E8) with company / report title text
#  - Populate column H with the summary "category" label for each account line
#  - Add the balancing check formula in H56 (=G56-F56)
#  - Restore the last active selection to D24

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header block (rows 2-8, column E) ------------------------------------
$ws.Range("E2").Value = "VSIG Pte. Ltd."
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = "Trial Balance"
$ws.Range("E8").Value = "December 2015"

# ---- Category labels (column H) -------------------------------------------
$ws.Range("H11").Value = "Bank Balances"
$ws.Range("H12").Value = "Bank Balances"
$ws.Range("H13").Value = "Bank Balances"
$ws.Range("H14").Value = "Trade Receivables"
$ws.Range("H15").Value = "Trade Receivables"
$ws.Range("H16").Value = "Plant and Equipment"
$ws.Range("H17").Value = "Plant and Equipment"
$ws.Range("H18").Value = "Plant and Equipment"
$ws.Range("H19").Value = "Plant and Equipment"
$ws.Range("H20").Value = "Deposits"
$ws.Range("H21").Value = "Prepayments"
$ws.Range("H22").Value = "Trade Payables"
$ws.Range("H23").Value = "Trade Payables"
$ws.Range("H26").Value = "GST Payables"
$ws.Range("H27").Value = "Accruals"
$ws.Range("H28").Value = "Amount owing to a Shareholder"
$ws.Range("H29").Value = "Income Tax Payables"
$ws.Range("H30").Value = "Share Capital"
$ws.Range("H31").Value = "Retained Profits"
$ws.Range("H32").Value = "Revenue"
$ws.Range("H33").Value = "Cost of Sales"
$ws.Range("H34").Value = "Accounting Fee"
$ws.Range("H35").Value = "Administrative Expenses"
$ws.Range("H36").Value = "Bank Charges"
$ws.Range("H37").Value = "Compilation Fee"
$ws.Range("H38").Value = "Depreciation"
$ws.Range("H39").Value = "Entertainment"
$ws.Range("H40").Value = "Freight Charges"
$ws.Range("H41").Value = "Internet Expenses"
$ws.Range("H42").Value = "Late Penalty"
$ws.Range("H43").Value = "Nominee Director Fee"
$ws.Range("H44").Value = "Office Supplies"
$ws.Range("H45").Value = "Postage and Courier"
$ws.Range("H46").Value = "Professional Fee"
$ws.Range("H47").Value = "Secretarial Fee"
$ws.Range("H48").Value = "Taxation Fee"
$ws.Range("H49").Value = "Telephone Expenses"
$ws.Range("H50").Value = "Salaries"
$ws.Range("H51").Value = "Skill Development Levy & SINDA"
$ws.Range("H52").Value = "Exchange Gain - Trade"
$ws.Range("H53").Value = "Exchange Gain - Non-trade"
$ws.Range("H54").Value = "Income Tax Expense"

# ---- Balancing formula ------------------------------------------------------
$ws.Range("H56").Formula = "=G56-F56"

# ---- Restore last selection --------------------------------------------------
$ws.Range("D24").Select()
